# Weekly price update: insert a new reading at the top of the Repollo table
# (row 147), pushing all subsequent rows down by one. This mirrors how the
# source feed prepends the newest weekly record to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 147; existing rows 147:242 shift to 148:243.
$ws.Rows(147).Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value = 44762
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = 100112006
$ws.Range("G147").Value = "Repollo"
$ws.Range("H147").Value = "Crespo record"
$ws.Range("I147").Value = "Segunda"
$ws.Range("J147").Value = 120
$ws.Range("K147").Value = 1000
$ws.Range("L147").Value = 1000
$ws.Range("M147").Value = 1000
$ws.Range("N147").Value = "$/unidad"
$ws.Range("O147").Value = "Provincia de Diguillín"
$ws.Range("P147").Value = 1000
$ws.Range("Q147").Value = 1
$ws.Range("R147").Value = "Hortaliza"
